# Add the newest S&P 500 earnings-growth data point.
# A new row is inserted at the top of the data table (row 2), pushing all
# existing rows down by one. The new row holds the latest date (2025-09-30,
# Excel serial 45930) and its corresponding earnings-growth value (11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (first data row),
# shifting rows 2:143 down to 3:144.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the latest data point.
$ws.Range("A2").Value = 45930
$ws.Range("B2").Value = 11

# The freshly inserted row picks up formatting from the row above (the bold
# header) by default; restore the plain data-row formatting by copying it
# from the row just below (which still carries the original date/number
# styles) using a formats-only paste so values are left untouched.
$ws.Range("A3:B3").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the saved selection state.
$ws.Range("B3").Select() | Out-Null
